$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G (Recorded By) to fit full names
# (49.166666666666664 compensates for the engine's internal padding so the
# stored column width ends up as exactly 50)
$ws.Columns.Item(7).ColumnWidth = 49.166666666666664

# Populate "Recorded By" (column G) values for recorded sessions
$ws.Range("G2").Value = "Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Eman Tantawi"
$ws.Range("G3").Value = "Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G4").Value = "Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G5").Value = "Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Nesma, Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G6").Value = "Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G7").Value = "Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G8").Value = "Dr. Majorelle Magdy, Dr. Asmaa Reda, Administrator, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki"
$ws.Range("G9").Value = "Dr. Majorelle Magdy, Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Manar Montaser"
$ws.Range("G10").Value = "Dr. Sara Wael, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad"
$ws.Range("G11").Value = "Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G12").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G13").Value = "D Wessam Atef, Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Range("G14").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G15").Value = "D Wessam Atef, Dr. Amal Awwad"
$ws.Range("G16").Value = "Dr. Nourhan Mohammad, Dr. Amal Awwad"
$ws.Range("G17").Value = "Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Nourhan Osama, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed"
$ws.Range("G18").Value = "Dr. Mohammad Safwat"
$ws.Range("G19").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G20").Value = "Dr. Nourham Mostafa"
$ws.Range("G22").Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Range("G23").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G24").Value = "Dr. Neveen Nashaat, Dr. Aya Emad, Dr. Monica, Dr. Yasmin, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Remon"
$ws.Range("G25").Value = "Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Marina Atef, Dr. Remon"
$ws.Range("G26").Value = "Dr. Gehad Salah, Dr. Youstina Magdy"
$ws.Range("G27").Value = "Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Remon"
$ws.Range("G28").Value = "Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Wafaa Ebida, Dr. Remon"
$ws.Range("G29").Value = "Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Remon, Dr. Naema Gomaa"
$ws.Range("G30").Value = "Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G31").Value = "Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G32").Value = "Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G33").Value = "Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Nesma, Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G34").Value = "Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G35").Value = "Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G36").Value = "Dr. Majorelle Magdy, Dr. Asmaa Reda, Administrator, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki"
$ws.Range("G37").Value = "Dr. Majorelle Magdy, Dr. Menna tu’Allah Medhat, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Manar Montaser"
$ws.Range("G38").Value = "Dr. Sara Wael, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad"
$ws.Range("G39").Value = "Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G40").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G41").Value = "D Wessam Atef, Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Range("G42").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G43").Value = "D Wessam Atef, Dr. Amal Awwad"
$ws.Range("G44").Value = "Dr. Nourhan Mohammad, Dr. Amal Awwad"
$ws.Range("G45").Value = "Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Nourhan Osama, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed"
$ws.Range("G46").Value = "Dr. Mohammad Safwat"
$ws.Range("G47").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G48").Value = "Dr. Nourham Mostafa"
$ws.Range("G50").Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Range("G51").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G52").Value = "Dr. Neveen Nashaat, Dr. Aya Emad, Dr. Monica, Dr. Yasmin, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Remon"
$ws.Range("G53").Value = "Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Marina Atef, Dr. Remon"
$ws.Range("G54").Value = "Dr. Gehad Salah, Dr. Youstina Magdy"
$ws.Range("G55").Value = "Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Remon"
$ws.Range("G56").Value = "Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Wafaa Ebida, Dr. Remon"
$ws.Range("G57").Value = "Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Remon, Dr. Naema Gomaa"
